$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells and match the existing header formatting
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every data row (rows 2 through 49)
for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = 67
    $ws.Cells.Item($row, 31).Value = 95
    $ws.Cells.Item($row, 32).Value = 0
}
